# Auto-generated Excel COM-interop script to apply value updates
# described in the commit "Model training severe changes".
# Updates numeric prediction values in columns A:C (rows 2-72) of the
# active worksheet, leaving headers (row 1) and all formatting untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object "object[,]" 71,3

$data[0,0] = 37.56304128944311
$data[0,1] = 32.11715364383557
$data[0,2] = 42.07352515068488
$data[1,0] = 30.36987052300136
$data[1,1] = 28.57669621917803
$data[1,2] = 32.08004942465744
$data[2,0] = 30.28320883254955
$data[2,1] = 28.3451558356164
$data[2,2] = 32.11361326027391
$data[3,0] = 36.99410115286576
$data[3,1] = 34.44114356164381
$data[3,2] = 39.36474706849311
$data[4,0] = 39.16227929394123
$data[4,1] = 35.17983572602729
$data[4,2] = 42.74264186301368
$data[5,0] = 32.0363216667348
$data[5,1] = 28.35816613698627
$data[5,2] = 35.11472328767122
$data[6,0] = 27.47084941484352
$data[6,1] = 22.29652372602735
$data[6,2] = 31.15838334246568
$data[7,0] = 36.21657615048041
$data[7,1] = 33.12124142465746
$data[7,2] = 39.07976021917801
$data[8,0] = 30.08999138608322
$data[8,1] = 24.07460964383555
$data[8,2] = 34.18231605479445
$data[9,0] = 36.79637967259583
$data[9,1] = 30.53580087671227
$data[9,2] = 41.56254465753423
$data[10,0] = 43.76599467157354
$data[10,1] = 38.40078761643834
$data[10,2] = 48.56167747945197
$data[11,0] = 38.08646146854762
$data[11,1] = 30.34447331506848
$data[11,2] = 44.28607671232874
$data[12,0] = 36.20138445961961
$data[12,1] = 29.52346641095885
$data[12,2] = 41.18519945205468
$data[13,0] = 42.93926871941659
$data[13,1] = 39.3509731506849
$data[13,2] = 46.32031956164374
$data[14,0] = 39.142818505827
$data[14,1] = 35.09621523287661
$data[14,2] = 42.77007353424656
$data[15,0] = 43.63463079288479
$data[15,1] = 38.6511619726027
$data[15,2] = 48.23861720547938
$data[16,0] = 30.31836509834386
$data[16,1] = 28.48670816438354
$data[16,2] = 32.21231813698626
$data[17,0] = 41.88252955005785
$data[17,1] = 34.13340635616434
$data[17,2] = 48.08768438356152
$data[18,0] = 44.6162161575682
$data[18,1] = 40.10430783561642
$data[18,2] = 48.83132558904105
$data[19,0] = 42.43928982921003
$data[19,1] = 38.77055408219173
$data[19,2] = 45.79736953424652
$data[20,0] = 33.00854988673064
$data[20,1] = 26.4028152328767
$data[20,2] = 38.15221775342455
$data[21,0] = 45.78763098698283
$data[21,1] = 41.64478794520539
$data[21,2] = 49.9160224657534
$data[22,0] = 39.15162879520206
$data[22,1] = 35.09889271232866
$data[22,2] = 42.77087243835614
$data[23,0] = 33.11637412826276
$data[23,1] = 30.50114421917804
$data[23,2] = 35.4084942465753
$data[24,0] = 42.6744052515504
$data[24,1] = 38.87428142465749
$data[24,2] = 45.97480021917809
$data[25,0] = 43.14935079615618
$data[25,1] = 32.86801501369857
$data[25,2] = 51.13303342465749
$data[26,0] = 34.9328945317249
$data[26,1] = 30.89325304109582
$data[26,2] = 38.62350290410945
$data[27,0] = 48.57379132992568
$data[27,1] = 44.38861917808209
$data[27,2] = 53.11317621917806
$data[28,0] = 40.96985111047498
$data[28,1] = 35.90994575342452
$data[28,2] = 45.57575813698629
$data[29,0] = 33.75950747141002
$data[29,1] = 28.2785675616438
$data[29,2] = 38.55127221917793
$data[30,0] = 30.28050722715186
$data[30,1] = 24.38587550684927
$data[30,2] = 34.53141610958904
$data[31,0] = 36.20762052858982
$data[31,1] = 29.26436821917803
$data[31,2] = 41.45359978082178
$data[32,0] = 44.81506114032575
$data[32,1] = 37.12384931506844
$data[32,2] = 51.48673326027392
$data[33,0] = 41.71100182430304
$data[33,1] = 37.72018060273972
$data[33,2] = 45.33229819178073
$data[34,0] = 31.96467435153001
$data[34,1] = 28.43794739726023
$data[34,2] = 35.24627002739723
$data[35,0] = 45.50130494322902
$data[35,1] = 41.30420591780812
$data[35,2] = 49.65190663013693
$data[36,0] = 40.94360446752533
$data[36,1] = 33.30478246575339
$data[36,2] = 47.22014191780815
$data[37,0] = 46.26958336754578
$data[37,1] = 40.58491243835611
$data[37,2] = 52.00030071232878
$data[38,0] = 36.77379066857483
$data[38,1] = 33.91590432876706
$data[38,2] = 39.3744484383561
$data[39,0] = 28.9842249683091
$data[39,1] = 25.79612876712326
$data[39,2] = 31.495992
$data[40,0] = 39.32503115436511
$data[40,1] = 33.00649972602739
$data[40,2] = 44.82599134246574
$data[41,0] = 40.9661333851291
$data[41,1] = 35.91302158904095
$data[41,2] = 45.5734072328767
$data[42,0] = 28.40893171294208
$data[42,1] = 24.8797594520547
$data[42,2] = 31.44568821917805
$data[43,0] = 44.51430322033659
$data[43,1] = 36.07128602739721
$data[43,2] = 51.03455846575334
$data[44,0] = 30.82816501710622
$data[44,1] = 25.68411210958903
$data[44,2] = 34.48812394520544
$data[45,0] = 35.60537002657935
$data[45,1] = 28.19327616438353
$data[45,2] = 41.15177599999996
$data[46,0] = 43.73690427806163
$data[46,1] = 38.39168767123286
$data[46,2] = 48.56750224657526
$data[47,0] = 33.43695984025077
$data[47,1] = 31.08203616438359
$data[47,2] = 35.63702279452048
$data[48,0] = 45.43307539398893
$data[48,1] = 41.21732065753415
$data[48,2] = 49.55390827397254
$data[49,0] = 30.3753930604511
$data[49,1] = 28.56044175342461
$data[49,2] = 32.06223353424647
$data[50,0] = 27.90763703564365
$data[50,1] = 23.78172032876708
$data[50,2] = 31.25129906849314
$data[51,0] = 46.57995674504183
$data[51,1] = 40.72196745205479
$data[51,2] = 51.77609501369862
$data[52,0] = 43.00318512097046
$data[52,1] = 39.41326882191778
$data[52,2] = 46.4283373150684
$data[53,0] = 42.54377718476108
$data[53,1] = 38.92118334246573
$data[53,2] = 46.24362586301367
$data[54,0] = 33.64603444476246
$data[54,1] = 28.02655167123286
$data[54,2] = 38.04612756164372
$data[55,0] = 30.28924247256861
$data[55,1] = 24.46318816438351
$data[55,2] = 34.41601523287671
$data[56,0] = 45.53232098630132
$data[56,1] = 41.33108843835606
$data[56,2] = 49.63432471232873
$data[57,0] = 34.06499597083069
$data[57,1] = 28.73592131506847
$data[57,2] = 38.76474860273963
$data[58,0] = 45.68932970980708
$data[58,1] = 41.51344043835607
$data[58,2] = 49.81773128767119
$data[59,0] = 38.71496092360111
$data[59,1] = 31.69943386301366
$data[59,2] = 44.7121695342465
$data[60,0] = 36.79637967259583
$data[60,1] = 30.53580087671227
$data[60,2] = 41.56254465753423
$data[61,0] = 43.70055842486193
$data[61,1] = 34.21935539726022
$data[61,2] = 51.4612361643834
$data[62,0] = 45.7033376165746
$data[62,1] = 41.5422271780821
$data[62,2] = 49.84310926027392
$data[63,0] = 27.43408228664888
$data[63,1] = 22.40959298630138
$data[63,2] = 30.84290104109581
$data[64,0] = 42.71383984951941
$data[64,1] = 36.03078739726016
$data[64,2] = 48.45456010958897
$data[65,0] = 40.00936840673341
$data[65,1] = 34.44760449315063
$data[65,2] = 45.13529315068484
$data[66,0] = 39.62826058420222
$data[66,1] = 36.39744810958894
$data[66,2] = 43.21659769863006
$data[67,0] = 35.44327503196344
$data[67,1] = 27.66789764383558
$data[67,2] = 41.30260493150685
$data[68,0] = 36.66425107176436
$data[68,1] = 33.85814224657532
$data[68,2] = 39.31424065753421
$data[69,0] = 29.50806388441347
$data[69,1] = 27.16761106849316
$data[69,2] = 31.54658805479446
$data[70,0] = 35.536729207115
$data[70,1] = 31.67840701369857
$data[70,2] = 38.79207123287664

$range = $ws.Range("A2:C72")
$range.Value = $data
